# Auto-generated edit script applying scheduled-runner price refresh to Seraph_Profits sheets.
# Source: unified diff of the workbook's canonical OOXML (commit: "chore: update Sheets via scheduled runner").
# All target cells are plain numeric literals (no formulas in this workbook), so values are written directly.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6420.4
$ws.Range("I62").Value = 5051
$ws.Range("J62").Value = 7333.3335
$ws.Range("K62").Value = 5051
$ws.Range("L62").Value = 7333.3335
$ws.Range("M62").Value = -4427
$ws.Range("N62").Value = -8581.333500000001
$ws.Range("H65").Value = 6420.4
$ws.Range("I65").Value = 5051
$ws.Range("J65").Value = 7333.3335
$ws.Range("K65").Value = 25255
$ws.Range("L65").Value = 36666.6675
$ws.Range("M65").Value = -22135
$ws.Range("N65").Value = -42906.6675
$ws.Range("H69").Value = 9578.429
$ws.Range("I69").Value = 9049.666999999999
$ws.Range("J69").Value = 9975
$ws.Range("K69").Value = 27149.001
$ws.Range("L69").Value = 29925
$ws.Range("M69").Value = -26275.001
$ws.Range("N69").Value = -31673
$ws.Range("H72").Value = 9578.429
$ws.Range("I72").Value = 9049.666999999999
$ws.Range("J72").Value = 9975
$ws.Range("K72").Value = 81447.003
$ws.Range("L72").Value = 89775
$ws.Range("M72").Value = -77079.003
$ws.Range("N72").Value = -98511
$ws.Range("H86").Value = 5999.25
$ws.Range("J86").Value = 5999.4287
$ws.Range("L86").Value = 5999.4287
$ws.Range("N86").Value = -8245.4287
$ws.Range("H88").Value = 5688.25
$ws.Range("J88").Value = 6792.6665
$ws.Range("L88").Value = 6792.6665
$ws.Range("N88").Value = -7604.6665
$ws.Range("H89").Value = 5999.25
$ws.Range("J89").Value = 5999.4287
$ws.Range("L89").Value = 29997.1435
$ws.Range("N89").Value = -41229.14350000001
$ws.Range("H91").Value = 5688.25
$ws.Range("J91").Value = 6792.6665
$ws.Range("L91").Value = 6792.6665
$ws.Range("N91").Value = -9600.666499999999
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2025
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 2025
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 1229
$ws.Range("N113").Value = -8408
$ws.Range("H116").Value = 4837.5
$ws.Range("I116").Value = 4658.3335
$ws.Range("K116").Value = 4658.3335
$ws.Range("M116").Value = -1216.3335
$ws.Range("H137").Value = 3385.7222
$ws.Range("I137").Value = 1562.25
$ws.Range("J137").Value = 4844.5
$ws.Range("K137").Value = 4686.75
$ws.Range("L137").Value = 14533.5
$ws.Range("M137").Value = -2136.75
$ws.Range("N137").Value = -19633.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1420.1
$ws.Range("I74").Value = 1134.1077
$ws.Range("K74").Value = 1134.1077
$ws.Range("M74").Value = -260.1077
$ws.Range("H77").Value = 1420.1
$ws.Range("I77").Value = 1134.1077
$ws.Range("K77").Value = 5670.538500000001
$ws.Range("M77").Value = -1302.538500000001
$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2500
$ws.Range("K102").Value = 2500
$ws.Range("M102").Value = -878
$ws.Range("H132").Value = 1442.4783
$ws.Range("I132").Value = 1442.4783
$ws.Range("K132").Value = 4327.4349
$ws.Range("M132").Value = -1797.4349

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2169.8
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2169.8
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H94").Value = 756.5172
$ws.Range("I94").Value = 794.2222
$ws.Range("K94").Value = 794.2222
$ws.Range("M94").Value = -343.2222

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 554
$ws.Range("I22").Value = 237
$ws.Range("K22").Value = 237
$ws.Range("M22").Value = 113
$ws.Range("H58").Value = 3720.318
$ws.Range("I58").Value = 1819.091
$ws.Range("K58").Value = 1819.091
$ws.Range("M58").Value = -1616.091
$ws.Range("H99").Value = 11461.639
$ws.Range("I99").Value = 8110.5
$ws.Range("J99").Value = 14812.777
$ws.Range("K99").Value = 8110.5
$ws.Range("L99").Value = 14812.777
$ws.Range("M99").Value = -6612.5
$ws.Range("N99").Value = -17808.777
$ws.Range("H122").Value = 2822.05
$ws.Range("I122").Value = 2682.125
$ws.Range("J122").Value = 3381.75
$ws.Range("K122").Value = 8046.375
$ws.Range("L122").Value = 10145.25
$ws.Range("M122").Value = -5596.375
$ws.Range("N122").Value = -15045.25
$ws.Range("H126").Value = 11461.639
$ws.Range("I126").Value = 8110.5
$ws.Range("J126").Value = 14812.777
$ws.Range("K126").Value = 24331.5
$ws.Range("L126").Value = 44438.331
$ws.Range("M126").Value = -21861.5
$ws.Range("N126").Value = -49378.331
$ws.Range("H132").Value = 2947.3333
$ws.Range("I132").Value = 2495.25
$ws.Range("K132").Value = 7485.75
$ws.Range("M132").Value = -4955.75
$ws.Range("H134").Value = 2770.5483
$ws.Range("I134").Value = 2100.611
$ws.Range("K134").Value = 6301.833
$ws.Range("M134").Value = -3766.833
$ws.Range("H136").Value = 3720.318
$ws.Range("I136").Value = 1819.091
$ws.Range("K136").Value = 5457.272999999999
$ws.Range("M136").Value = -2907.272999999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 870.8333
$ws.Range("I117").Value = 246.5
$ws.Range("K117").Value = 739.5
$ws.Range("M117").Value = 2702.5
$ws.Range("H128").Value = 3979891
$ws.Range("I128").Value = 3979891
$ws.Range("K128").Value = 11939673
$ws.Range("M128").Value = -11934693

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4597.4165
$ws.Range("I132").Value = 1874.1666
$ws.Range("K132").Value = 5622.4998
$ws.Range("M132").Value = -3092.4998

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 49999
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 49999
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H82").Value = 3000.3928
$ws.Range("I82").Value = 3187.9524
$ws.Range("J82").Value = 2437.7144
$ws.Range("K82").Value = 3187.9524
$ws.Range("L82").Value = 2437.7144
$ws.Range("M82").Value = -2826.9524
$ws.Range("N82").Value = -3159.7144
$ws.Range("H85").Value = 3000.3928
$ws.Range("I85").Value = 3187.9524
$ws.Range("J85").Value = 2437.7144
$ws.Range("K85").Value = 3187.9524
$ws.Range("L85").Value = 2437.7144
$ws.Range("M85").Value = -1939.9524
$ws.Range("N85").Value = -4933.7144
$ws.Range("H104").Value = 49999
$ws.Range("J104").Value = 49999
$ws.Range("L104").Value = 49999
$ws.Range("N104").Value = -56987

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 514000
$ws.Range("J49").Value = 28000
$ws.Range("L49").Value = 28000
$ws.Range("N49").Value = -28460
$ws.Range("H100").Value = 2149.3333
$ws.Range("I100").Value = 1880
$ws.Range("K100").Value = 3760
$ws.Range("M100").Value = -3219
